# Insert a new data row at row 248 (pushing the existing rows 248:346
# down to 249:347) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("248:248").Insert()

$ws.Cells.Item(248, 1).Value = 6
$ws.Cells.Item(248, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(248, 3).Value = "Metropolitana"
$ws.Cells.Item(248, 4).Value = 45027
$ws.Cells.Item(248, 5).Value = 13
$ws.Cells.Item(248, 6).Value = "Fruta"
$ws.Cells.Item(248, 7).Value = 100101
$ws.Cells.Item(248, 8).Value = "Berries"
$ws.Cells.Item(248, 9).Value = 100101004
$ws.Cells.Item(248, 10).Value = "Frambuesa"
$ws.Cells.Item(248, 11).Value = "Sin especificar"
$ws.Cells.Item(248, 12).Value = "Primera"
$ws.Cells.Item(248, 13).Value = 50
$ws.Cells.Item(248, 14).Value = 9000
$ws.Cells.Item(248, 15).Value = 9000
$ws.Cells.Item(248, 16).Value = 9000
$ws.Cells.Item(248, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(248, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(248, 19).Value = 4500
$ws.Cells.Item(248, 20).Value = 2
